# The commit marks task #4 ("Initialize React.js project structure") as
# completed by ticking its "Completed" checkbox (cell E8 on the
# "To Do List Check Box" sheet). The checkbox form control's linked cell
# is $E$8, so setting that cell to TRUE is equivalent to the user clicking
# the checkbox. Every downstream figure (Completed Score in H10, the
# % Completed in H12, and the chart/label that mirror H12) recalculates
# automatically from this single input change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do List Check Box")

$ws.Range("E8").Value = $true

$wb.Save()
